# Updates cryptos list: Price (D) and Volume(1h) (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.393.19'
$ws.Range("E2").Value = '  +2.42%  '
$ws.Range("D3").Value = '2.540.03'
$ws.Range("E3").Value = '  +2.88%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''541.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.05%  '
$ws.Range("D6").Value = '''145.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.79%  '
$ws.Range("E7").Value = '  -0.44%  '
$ws.Range("E8").Value = '  +0.74%  '
$ws.Range("D9").Value = '2.573.05'
$ws.Range("E9").Value = '  +2.99%  '
$ws.Range("E10").Value = '  +1.75%  '
$ws.Range("E11").Value = '  +1.64%  '
$ws.Range("E12").Value = '  -1.13%  '
$ws.Range("E13").Value = '  +2.60%  '
$ws.Range("D14").Value = '2.986.71'
$ws.Range("E14").Value = '  +2.19%  '
$ws.Range("D15").Value = '''24.33'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.57%  '
$ws.Range("D16").Value = '60.349.39'
$ws.Range("E16").Value = '  +2.41%  '
$ws.Range("E17").Value = '  +3.90%  '
$ws.Range("D18").Value = '2.562.92'
$ws.Range("E18").Value = '  +2.61%  '
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("E20").Value = '  +0.90%  '
$ws.Range("D21").Value = '''329.15'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.47%  '
$ws.Range("E22").Value = '  +0.38%  '
$ws.Range("D23").Value = '''5.94'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.06%  '
$ws.Range("D24").Value = '''63.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.64%  '
$ws.Range("E25").Value = '  -0.12%  '
$ws.Range("E26").Value = '  +3.95%  '
$ws.Range("E27").Value = '  -0.71%  '
$ws.Range("D28").Value = '''8.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.03%  '
$ws.Range("D29").Value = '''7.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.86%  '
$ws.Range("D30").Value = '0.0₃0802'
$ws.Range("E30").Value = '  +3.10%  '
$ws.Range("E31").Value = '  +0.98%  '
$ws.Range("E32").Value = '  -3.83%  '
$ws.Range("D33").Value = '''162.54'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.75%  '
$ws.Range("E34").Value = '  +6.07%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Value = '''18.85'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.64%  '
$ws.Range("D37").Value = '''4.53'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.13%  '
$ws.Range("E38").Value = '  +0.83%  '
$ws.Range("D39").Value = '''5.74'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.41%  '
$ws.Range("D40").Value = '''37.27'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.62%  '
$ws.Range("D41").Value = '''305.91'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.57%  '
$ws.Range("D42").Value = '''0.847'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.80%  '
$ws.Range("E43").Value = '  +1.11%  '
$ws.Range("D44").Value = '''0.610'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.55%  '
$ws.Range("D45").Value = '''0.992'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.47%  '
$ws.Range("E46").Value = '  +0.95%  '
$ws.Range("D47").Value = '''19.16'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.44%  '
$ws.Range("E48").Value = '  +1.31%  '
$ws.Range("D49").Value = '''125.24'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.88%  '
$ws.Range("D50").Value = '''0.0527'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.60%  '
$ws.Range("E51").Value = '  +0.76%  '
